$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shp = $master.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$found = $tr.Replace("4/12/16", "4/14/2017")
Write-Output "done"
